# Update 2021 HWL2 First Batch
#
# Extends the Poland-Lithuania exchange-rate workbook with year columns
# 2016-2050 (35 new years) so future data batches have somewhere to land.
#
#  - "Data Clio Infra Format": the wide/pivoted sheet. Year headers run
#    across row 1 (I1=1500 ... TD1=2015); we append TE1:UM1 = 2016..2050
#    as text headers, mirroring the existing year-header cells. Row 2
#    (the single data row) gets no values for these new years - they stay
#    blank, same as all the other not-yet-observed years.
#
#  - "Data Long Format": the tall/long sheet. Columns E:F used to hold the
#    "year"/"value" pair. We insert 35 blank columns at E so E1:AM1 can
#    hold the same 2016..2050 text headers as the other sheet; the
#    pre-existing "year"/"value" columns (and their data, rows 2-10) slide
#    right to AN:AO untouched.
#
#  - "Metadata" sheet is untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Data Clio Infra Format" - append year columns TE:UM (2016-2050)
# ---------------------------------------------------------------------
$wsWide = $wb.Worksheets.Item("Data Clio Infra Format")

$firstNewCol = 525   # column TE (one past TD, which holds 2015)
$yearCount = 35      # 2016 .. 2050 inclusive

$wideHeaderRange = $wsWide.Range($wsWide.Cells.Item(1, $firstNewCol), $wsWide.Cells.Item(1, $firstNewCol + $yearCount - 1))
# Force text storage (matches how the existing "1500".."2015" year headers
# are stored as text, not numbers) before writing the values.
$wideHeaderRange.NumberFormat = "@"

for ($i = 0; $i -lt $yearCount; $i++) {
    $col = $firstNewCol + $i
    $year = 2016 + $i
    $wsWide.Cells.Item(1, $col).Value = "$year"
    # Row 2 (the only data row) is left blank for these new, not-yet-populated years.
}

# ---------------------------------------------------------------------
# Sheet 2: "Data Long Format" - insert year columns E:AM (2016-2050),
# pushing the old "year"/"value" columns from E:F to AN:AO.
# ---------------------------------------------------------------------
$wsLong = $wb.Worksheets.Item("Data Long Format")

$insertAtCol = 5          # column E
$longYearCount = 35       # 2016 .. 2050 inclusive, lands in E:AM

$insertRange = $wsLong.Range($wsLong.Cells.Item(1, $insertAtCol), $wsLong.Cells.Item(1, $insertAtCol + $longYearCount - 1))
$insertRange.EntireColumn.Insert()

$longHeaderRange = $wsLong.Range($wsLong.Cells.Item(1, $insertAtCol), $wsLong.Cells.Item(1, $insertAtCol + $longYearCount - 1))
$longHeaderRange.NumberFormat = "@"

for ($i = 0; $i -lt $longYearCount; $i++) {
    $col = $insertAtCol + $i
    $year = 2016 + $i
    $wsLong.Cells.Item(1, $col).Value = "$year"
    # Rows 2-10 under these new columns stay blank - no long-format
    # observations exist yet for 2016-2050.
}

Write-Host "Added year columns 2016-2050 to 'Data Clio Infra Format' and 'Data Long Format'."
